$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        $ParaRange,
        [string]$InnerXml
    )
    $pkg = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
$InnerXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    [void]$ParaRange.InsertXML($pkg)
}

# --- 1) Append the brand-new paragraph 11 at the very end of the document ---
# (Paragraph.Range on the LAST paragraph includes the document's terminal
#  mark, so we collapse a zero-length range at its End instead of handing
#  InsertXML the whole paragraph range.)
$lastParaRange = $d.Paragraphs($d.Paragraphs.Count).Range
$endPoint = $d.Range($lastParaRange.End, $lastParaRange.End)
Set-ParagraphXml $endPoint '<w:p><w:r><w:t>11 – Velocidad 5 – 10 - 5</w:t></w:r></w:p>'

# --- 2) Paragraph 1: "1 – Salto Abalacob" -> "1 – Salto Abalakov" split into
#        several runs, flanked by proofErr spell-check marks ---
$p1 = $d.Paragraphs(1).Range
$p1Xml = @'
<w:p>
  <w:r><w:t xml:space="preserve">1 – Salto </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Abala</w:t></w:r>
  <w:r><w:t>k</w:t></w:r>
  <w:r><w:t>o</w:t></w:r>
  <w:r><w:t>v</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@
Set-ParagraphXml $p1 $p1Xml

# --- 3) Paragraph 2: "2 – Salto cmj" split into runs with proofErr marks ---
$p2 = $d.Paragraphs(2).Range
$p2Xml = @'
<w:p>
  <w:r><w:t xml:space="preserve">2 – Salto </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>cmj</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@
Set-ParagraphXml $p2 $p2Xml

# --- 4) Paragraph 3: "3 – Salto sj" split into runs with proofErr marks ---
$p3 = $d.Paragraphs(3).Range
$p3Xml = @'
<w:p>
  <w:r><w:t xml:space="preserve">3 – Salto </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>sj</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@
Set-ParagraphXml $p3 $p3Xml

# --- 5) Paragraph 6: "6 – Velocidad 10 mts" split into runs with proofErr ---
$p6 = $d.Paragraphs(6).Range
$p6Xml = @'
<w:p>
  <w:r><w:t xml:space="preserve">6 – Velocidad 10 </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>mts</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@
Set-ParagraphXml $p6 $p6Xml

# --- 6) Paragraph 9: "9 – Sentadilla Bulgara" -> "9 – Sentadilla Búlgara",
#        split into two runs; the _GoBack bookmark moves here (off of
#        paragraph 10) ---
$p9 = $d.Paragraphs(9).Range
$p9Xml = @'
<w:p>
  <w:r><w:t xml:space="preserve">9 – Sentadilla </w:t></w:r>
  <w:r><w:t>Búlgara</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@
Set-ParagraphXml $p9 $p9Xml

# --- 7) Paragraph 10: same text, but the _GoBack bookmark is removed
#        (it now lives on paragraph 9) ---
$p10 = $d.Paragraphs(10).Range
$p10Xml = '<w:p><w:r><w:t>10 – Peso muerto 1 pierna</w:t></w:r></w:p>'
Set-ParagraphXml $p10 $p10Xml
